$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new "Date" column before the existing "Problem Statement"
#    column (old B -> C, old C -> D, old D -> E, old E -> F).
# ---------------------------------------------------------------------------
$ws.Columns("B:B").Insert()

# ---------------------------------------------------------------------------
# 2. Header row (row 1)
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "Date"
$ws.Range("B1").Font.Bold = $true

$ws.Range("G1").Value = "who"
$ws.Range("G1").Font.Bold = $true

# ---------------------------------------------------------------------------
# 3. Wrap text for the "Problem Statement" column (now column C).
#    Apply to the body first, then the header, so the style table ends up
#    in the same order as a real Excel edit (plain+wrap before bold+wrap).
#    Only the rows that actually hold data (3-7) are touched so blank rows
#    below stay untouched, same as the target workbook.
# ---------------------------------------------------------------------------
$ws.Range("C3:C7").WrapText = $true
$ws.Range("C1").WrapText = $true

# Column widths
$ws.Columns("B:B").ColumnWidth = 9.5
$ws.Columns("C:C").ColumnWidth = 38.85

# ---------------------------------------------------------------------------
# 4. Row 3 data
# ---------------------------------------------------------------------------
$ws.Range("B3").Formula = "=DATE(2024,12,21)"
$ws.Range("B3").NumberFormat = "mm-dd-yy"

$ws.Range("C3").Value = "Remove Duplicates from sorted array"
$ws.Range("F3").Value = "used haspmap or hashset, also solved without extra space"

# ---------------------------------------------------------------------------
# 5. Row 4 data
# ---------------------------------------------------------------------------
$ws.Range("B4").Formula = "=DATE(2024,12,21)"
$ws.Range("B4").NumberFormat = "mm-dd-yy"

$ws.Range("C4").Value = "find target in rotated sorted array .  Use O(logn)"

$ws.Range("F4").Value = "first check left part is sorted or right. Then in each part check target lies in left part or right part"
$ws.Range("F4").WrapText = $true

$ws.Range("G4").Value = "took help"

# ---------------------------------------------------------------------------
# 6. Row 5 data (contains the mixed-font explanation)
# ---------------------------------------------------------------------------
$ws.Range("B5").Formula = "=DATE(2024,12,21)"
$ws.Range("B5").NumberFormat = "mm-dd-yy"

$ws.Range("C5").Value = "find minimum in rotated sorted array in logn time complexity"

$ws.Range("F5").Value = "If arr[mid] > arr[high], the minimum lies in the right half, so set low = mid + 1.`nIf arr[mid] <= arr[high], the minimum lies in the left half or at mid, so set high = mid.`nWhen low == high, the pointer low (or high) will point to the minimum element."
$ws.Range("F5").WrapText = $true

$run1 = $ws.Range("F5").Characters(4, 20)
$run1.Font.Size = 10
$run1.Font.Name = "Arial Unicode MS"

$run2 = $ws.Range("F5").Characters(24, 45)
$run2.Font.Size = 11
$run2.Font.Name = "Calibri"

$run3 = $ws.Range("F5").Characters(69, 13)
$run3.Font.Size = 10
$run3.Font.Name = "Arial Unicode MS"

$run4 = $ws.Range("F5").Characters(82, 170)
$run4.Font.Size = 11
$run4.Font.Name = "Calibri"

$ws.Range("G5").Value = "took help"

# ---------------------------------------------------------------------------
# 7. Row 6 data
# ---------------------------------------------------------------------------
$ws.Range("B6").Formula = "=DATE(2024,12,21)"
$ws.Range("B6").NumberFormat = "mm-dd-yy"

$ws.Range("C6").Value = "find peak element. A peak element is an element that is strictly greater than its neighbors. Use O(logn)"

$ws.Range("F6").Value = "solved in O(n)`n  while (low < high) {`n            int mid = low + (high - low) / 2;`n            // Check if mid is a peak`n            if (arr[mid] > arr[mid + 1]) {`n                high = mid; // Peak is in the left half or at mid`n            } else {`n                low = mid + 1; // Peak is in the right half`n            }`n        }`n        return low; // or high, as low == high"
$ws.Range("F6").WrapText = $true

$ws.Range("G6").Value = "took help"

# ---------------------------------------------------------------------------
# 8. Row 7 data
# ---------------------------------------------------------------------------
$ws.Range("B7").Formula = "=DATE(2024,12,22)"
$ws.Range("B7").NumberFormat = "mm-dd-yy"

$ws.Range("C7").Value = "move zeroes to end"

# ---------------------------------------------------------------------------
# 9. Row heights for the wrapped rows
# ---------------------------------------------------------------------------
$ws.Rows(4).RowHeight = 28.8
$ws.Rows(5).RowHeight = 86.4
$ws.Rows(6).RowHeight = 201.6

# ---------------------------------------------------------------------------
# 10. View state
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("F7").Select()
